$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# Ensure column B keeps a text number format so date-like strings are not
# auto-converted into Excel date serials.
$ws1.Range("B2:B17").NumberFormat = "@"

$ws1.Range("B2").Value = "2024-12-29"
$ws1.Range("D2").Value = 97
$ws1.Range("J2").Value = $null

$ws1.Range("B3").Value = "2025-01-05"
$ws1.Range("D3").Value = 105
$ws1.Range("J3").Value = $null

$ws1.Range("B4").Value = "2025-01-12"
$ws1.Range("D4").Value = 90
$ws1.Range("J4").Value = $null

$ws1.Range("B5").Value = "2025-01-19"
$ws1.Range("D5").Value = 89
$ws1.Range("J5").Value = $null

$ws1.Range("B6").Value = "2025-01-26"
$ws1.Range("D6").Value = 89
$ws1.Range("J6").Value = $null

$ws1.Range("B7").Value = "2025-02-02"
$ws1.Range("D7").Value = 89
$ws1.Range("J7").Value = $null

$ws1.Range("B8").Value = "2025-02-09"
$ws1.Range("D8").Value = 90
$ws1.Range("J8").Value = $null

$ws1.Range("B9").Value = "2025-02-16"
$ws1.Range("D9").Value = 91
$ws1.Range("J9").Value = $null

$ws1.Range("B10").Value = "2025-02-23"
$ws1.Range("D10").Value = 89
$ws1.Range("J10").Value = $null

$ws1.Range("B11").Value = "2025-03-02"
$ws1.Range("D11").Value = 90
$ws1.Range("J11").Value = $null

$ws1.Range("B12").Value = "2025-03-09"
$ws1.Range("D12").Value = 91
$ws1.Range("J12").Value = $null

$ws1.Range("B13").Value = "2025-03-16"
$ws1.Range("D13").Value = 90
$ws1.Range("J13").Value = $null

$ws1.Range("B14").Value = "2025-03-23"
$ws1.Range("D14").Value = 89
$ws1.Range("J14").Value = $null

$ws1.Range("B15").Value = "2025-03-30"
$ws1.Range("D15").Value = 88
$ws1.Range("J15").Value = $null

$ws1.Range("B16").Value = "2025-04-06"
$ws1.Range("D16").Value = 88
$ws1.Range("J16").Value = $null

$ws1.Range("B17").Value = "2025-04-13"
$ws1.Range("D17").Value = 88
$ws1.Range("J17").Value = $null

# Update Summary sheet aggregate metrics.
# Keep these as text cells (matching the source workbook's inlineStr cells)
# so numeric-looking and date-looking strings are not auto-converted.
$ws2.Range("B9:B15").NumberFormat = "@"
$ws2.Range("B9").Value = "1453"
$ws2.Range("B10").Value = "740"
$ws2.Range("B11").Value = "381"
$ws2.Range("B12").Value = "105"
$ws2.Range("B13").Value = "2025-01-05"
$ws2.Range("B14").Value = "88"
$ws2.Range("B15").Value = "2025-04-13"

Write-Host "Edit applied successfully"
